$d = $word.ActiveDocument

# 1. Update the timestamp in the "Fecha" field result text: 08:11 p. m. -> 08:42 p. m.
$d.Content.Find.Execute("12/06/2022 08:11 p. m.", $true, $false, $false, $false, $false, $true, 1, $false, "12/06/2022 08:42 p. m.", 2) | Out-Null

# 2. Justify ("both") the narrative paragraphs and bullet items that follow the banner table.
#    These are Paragraphs 11-22 in document order (everything after the version/date table).
for ($i = 11; $i -le 22; $i++) {
    $d.Paragraphs.Item($i).Alignment = 3
}

# 3. Normalize runs that were needlessly split (no formatting differences between them) by
#    re-finding and replacing the text in place with itself; the engine collapses adjacent
#    same-format runs that participate in a replace into a single run, exactly like Word does.

# "registro" + "_maestro" -> "registro_maestro" (occurs in two bullets), which also merges
# each bullet's trailing split runs into one run as a side effect of the in-place replace.
$d.Content.Find.Execute("registro_maestro", $true, $false, $false, $false, $false, $true, 1, $false, "registro_maestro", 2) | Out-Null

# Collapse the fully-split "social" bullet into a single run.
$d.Content.Find.Execute("Se modificó el acceso de registrov2 en la página social, ahora redireccionando a registros.", $true, $false, $false, $false, $false, $true, 1, $false, "Se modificó el acceso de registrov2 en la página social, ahora redireccionando a registros.", 2) | Out-Null
